$wb = $excel.ActiveWorkbook

# --- Config sheet: RunInParallel value changes from "No" to "Yes" ---
$configWs = $wb.Worksheets.Item("Config")
$configWs.Range("B3").Value = "Yes"
$configWs.Range("B13").Select() | Out-Null

# --- Test Cases sheet: add two new failed-testcase rows, then update Execute value on row 2 ---
$testsWs = $wb.Worksheets.Item("Test Cases")

$testsWs.Range("C13").Value = "5"
$testsWs.Range("A13").Value = "1111"
$testsWs.Range("B13").Value = "Failed TestCase"
$testsWs.Range("D13").Value = "Sanity"

$testsWs.Range("A14").Value = "2222"
$testsWs.Range("B14").Value = "Failed TestCase2"
$testsWs.Range("C14").Value = "6"
$testsWs.Range("D14").Value = "Sanity"

$testsWs.Range("E2").Value = "TestCaseNumber=1111,2222"

$testsWs.Range("D10").Select() | Out-Null
